# Replace the placeholder narrative text "dsdsds" with "Stage_narrative"
# on slides 2 through 5 (the "Stage_narrative" textbox, shape named
# "TextBox 1", is identical on each journey-stage slide).

$p = $ppt.ActivePresentation

for ($i = 2; $i -le 5; $i++) {
    $s = $p.Slides.Item($i)
    foreach ($shp in $s.Shapes) {
        if ($shp.Name -eq "TextBox 1") {
            $shp.TextFrame.TextRange.Text = "Stage_narrative"
        }
    }
}
